$wb = $excel.ActiveWorkbook

# Auto-generated from the OOXML diff: each block updates the market-price /
# leve-profit columns (H..N) for a single row, refreshed by the scheduled
# price-scraper run. Two rows (LTW!N68, LTW!N71) lose their HQ-profit cell
# entirely because HQ price data (J/L) dropped to 0 for that item.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1884.5385; $ws.Range("I28").Value = 133.33333; $ws.Range("J28").Value = 5824.75; $ws.Range("K28").Value = 133.33333; $ws.Range("L28").Value = 5824.75; $ws.Range("M28").Value = 351.66667; $ws.Range("N28").Value = -6794.75
$ws.Range("H86").Value = 1840.1428; $ws.Range("I86").Value = 1798.6; $ws.Range("K86").Value = 1798.6; $ws.Range("M86").Value = -675.5999999999999
$ws.Range("H89").Value = 1840.1428; $ws.Range("I89").Value = 1798.6; $ws.Range("K89").Value = 8993; $ws.Range("M89").Value = -3377
$ws.Range("H94").Value = 2425; $ws.Range("I94").Value = 2425; $ws.Range("K94").Value = 2425; $ws.Range("M94").Value = -1974
$ws.Range("H107").Value = 1402.5; $ws.Range("I107").Value = 772.3333; $ws.Range("J107").Value = 1875.125; $ws.Range("K107").Value = 772.3333; $ws.Range("L107").Value = 1875.125; $ws.Range("M107").Value = 1147.6667; $ws.Range("N107").Value = -5715.125
$ws.Range("H132").Value = 1167.15; $ws.Range("I132").Value = 1076.9375; $ws.Range("J132").Value = 1528; $ws.Range("K132").Value = 3230.8125; $ws.Range("L132").Value = 4584; $ws.Range("M132").Value = -700.8125; $ws.Range("N132").Value = -9644
$ws.Range("H137").Value = 1488.9; $ws.Range("I137").Value = 814.8333; $ws.Range("J137").Value = 2500; $ws.Range("K137").Value = 2444.4999; $ws.Range("L137").Value = 7500; $ws.Range("M137").Value = 105.5001000000002; $ws.Range("N137").Value = -12600
$ws.Range("H141").Value = 8030.5; $ws.Range("I141").Value = 3300; $ws.Range("K141").Value = 9900; $ws.Range("M141").Value = -4720
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 420; $ws.Range("I110").Value = 420; $ws.Range("K110").Value = 420; $ws.Range("M110").Value = 1625
$ws.Range("H122").Value = 1773.4445; $ws.Range("I122").Value = 1107.6666; $ws.Range("K122").Value = 3322.9998; $ws.Range("M122").Value = -872.9998000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2492; $ws.Range("I105").Value = 2205.524; $ws.Range("K105").Value = 2205.524; $ws.Range("M105").Value = -458.5239999999999
$ws.Range("H107").Value = 1048.9333; $ws.Range("I107").Value = 866.4286; $ws.Range("J107").Value = 1208.625; $ws.Range("K107").Value = 866.4286; $ws.Range("L107").Value = 1208.625; $ws.Range("M107").Value = 1053.5714; $ws.Range("N107").Value = -5048.625
$ws.Range("H122").Value = 40000; $ws.Range("J122").Value = 40000; $ws.Range("L122").Value = 40000; $ws.Range("N122").Value = -49800
$ws.Range("H134").Value = 4440.0938; $ws.Range("I134").Value = 4556.759; $ws.Range("K134").Value = 13670.277; $ws.Range("M134").Value = -11135.277
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 483.83334; $ws.Range("I16").Value = 419.4; $ws.Range("K16").Value = 419.4; $ws.Range("M16").Value = -132.4
$ws.Range("H58").Value = 1674513.5; $ws.Range("I58").Value = 2900423; $ws.Range("J58").Value = 2818.6365; $ws.Range("K58").Value = 2900423; $ws.Range("L58").Value = 2818.6365; $ws.Range("M58").Value = -2900220; $ws.Range("N58").Value = -3224.6365
$ws.Range("H105").Value = 1123; $ws.Range("I105").Value = 1130.6666; $ws.Range("J105").Value = 1100; $ws.Range("K105").Value = 1130.6666; $ws.Range("L105").Value = 1100; $ws.Range("M105").Value = 616.3334; $ws.Range("N105").Value = -4594
$ws.Range("H113").Value = 483.83334; $ws.Range("I113").Value = 419.4; $ws.Range("K113").Value = 419.4; $ws.Range("M113").Value = 1750.6
$ws.Range("H136").Value = 1674513.5; $ws.Range("I136").Value = 2900423; $ws.Range("J136").Value = 2818.6365; $ws.Range("K136").Value = 8701269; $ws.Range("L136").Value = 8455.9095; $ws.Range("M136").Value = -8698719; $ws.Range("N136").Value = -13555.9095
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1980; $ws.Range("I64").Value = 490; $ws.Range("K64").Value = 1470; $ws.Range("M64").Value = -1200
$ws.Range("H67").Value = 1980; $ws.Range("I67").Value = 490; $ws.Range("K67").Value = 1470; $ws.Range("M67").Value = -534
$ws.Range("H131").Value = 15991.218; $ws.Range("J131").Value = 17446.809; $ws.Range("L131").Value = 52340.427; $ws.Range("N131").Value = -62420.427
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2932; $ws.Range("I80").Value = 3660; $ws.Range("J80").Value = 1840; $ws.Range("K80").Value = 3660; $ws.Range("L80").Value = 1840; $ws.Range("M80").Value = -2662; $ws.Range("N80").Value = -3836
$ws.Range("H83").Value = 2932; $ws.Range("I83").Value = 3660; $ws.Range("J83").Value = 1840; $ws.Range("K83").Value = 18300; $ws.Range("L83").Value = 9200; $ws.Range("M83").Value = -13308; $ws.Range("N83").Value = -19184
$ws.Range("H113").Value = 754.5333000000001; $ws.Range("I113").Value = 343.625; $ws.Range("K113").Value = 343.625; $ws.Range("M113").Value = 1826.375
$ws.Range("H132").Value = 1749496.8; $ws.Range("I132").Value = 1924121.5; $ws.Range("J132").Value = 3249.5; $ws.Range("K132").Value = 5772364.5; $ws.Range("L132").Value = 9748.5; $ws.Range("M132").Value = -5769834.5; $ws.Range("N132").Value = -14808.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1662.6666; $ws.Range("I68").Value = 1662.6666; $ws.Range("J68").Value = 0; $ws.Range("K68").Value = 1662.6666; $ws.Range("L68").Value = 0; $ws.Range("M68").Value = -913.6666
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1662.6666; $ws.Range("I71").Value = 1662.6666; $ws.Range("J71").Value = 0; $ws.Range("K71").Value = 8313.333000000001; $ws.Range("L71").Value = 0; $ws.Range("M71").Value = -4569.333000000001
$ws.Range("N71").ClearContents()
$ws.Range("H82").Value = 1759.0667; $ws.Range("J82").Value = 2187; $ws.Range("L82").Value = 2187; $ws.Range("N82").Value = -2909
$ws.Range("H85").Value = 1759.0667; $ws.Range("J85").Value = 2187; $ws.Range("L85").Value = 2187; $ws.Range("N85").Value = -4683
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4666.3335; $ws.Range("I62").Value = 3999.5; $ws.Range("K62").Value = 3999.5; $ws.Range("M62").Value = -3375.5
$ws.Range("H65").Value = 4666.3335; $ws.Range("I65").Value = 3999.5; $ws.Range("K65").Value = 19997.5; $ws.Range("M65").Value = -16877.5
$ws.Range("H107").Value = 768.2308; $ws.Range("I107").Value = 568.4; $ws.Range("J107").Value = 1434.3334; $ws.Range("K107").Value = 1705.2; $ws.Range("L107").Value = 4303.0002; $ws.Range("M107").Value = 214.8000000000002; $ws.Range("N107").Value = -8143.0002
$ws.Range("H113").Value = 521.3182; $ws.Range("I113").Value = 343.53845; $ws.Range("K113").Value = 1030.61535; $ws.Range("M113").Value = 1139.38465
$ws.Range("H122").Value = 39124.332; $ws.Range("I122").Value = 99698.5; $ws.Range("K122").Value = 299095.5; $ws.Range("M122").Value = -296645.5
